$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title heading.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Balloonies Farm Free Slot - Full Game
#    Review") right before the final ("Prompt: ...") paragraph. First split
#    off a fresh empty paragraph, then stamp its own (now non-boundary-
#    spanning) range with the desired run XML so neighboring paragraphs are
#    left untouched.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Balloonies Farm Free Slot - Full Game Review</w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml)

# 3. Replace the old "Prompt: ..." text in the final paragraph with the meta
#    description text, preserving the existing (italic) run formatting.
$d.Content.Find.Execute(
    "Prompt: Create a feature image that captures the fun and playful nature of Balloonies Farm, while also featuring a happy Maya warrior with glasses. The image should be in colorful cartoon style and prominently feature the various farm animals, balloons, and the Maya warrior. The Maya warrior should be depicted wearing glasses and looking joyful and excited as they play Balloonies Farm. They could be holding a handful of colorful balloons while surrounded by the fun and playful farm animals, such as the Blue Dog, Sheep, Cow, and Pig, all of which should be depicted as balloon-shaped. The background should feature a scenic farm landscape with rolling hills and clear skies, with the reels of the slot machine woven into the design. Overall, the image should be lighthearted and convey a sense of fun and enjoyment, inviting viewers to try out the Balloonies Farm slot game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play the fun and exciting Balloonies Farm slot for free. Read our full game review and discover the features, graphics, and extra chances to win.",
    2
)
